# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing header cell (AC1) onto the new header
# cells so the new headers look like the rest of the header row (bold, border,
# centered alignment, etc).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Fill in the team record for every player row (rows 2 through 43).
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 90   # AD = Wins
    $ws.Cells.Item($r, 31).Value2 = 72   # AE = Losses
    $ws.Cells.Item($r, 32).Value2 = 0    # AF = Ties
}

Write-Host "Team record columns added."
